$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 12:03"

# Update numeric data cells (Casos totales / Casos activos / Recuperados / Muertes)
# Row 4
$ws.Range("B4").Value = 62989
$ws.Range("C4").Value = 38002
$ws.Range("D4").Value = 16567
$ws.Range("E4").Value = 8420

# Row 5
$ws.Range("B5").Value = 50771
$ws.Range("C5").Value = 22553
$ws.Range("D5").Value = 22948
$ws.Range("E5").Value = 5270

# Row 6
$ws.Range("B6").Value = 17429
$ws.Range("C6").Value = 6935
$ws.Range("D6").Value = 8662
$ws.Range("E6").Value = 1832

# Row 7
$ws.Range("B7").Value = 16080
$ws.Range("C7").Value = 5824
$ws.Range("D7").Value = 7640
$ws.Range("E7").Value = 2616

# Row 9
$ws.Range("B9").Value = 12210
$ws.Range("C9").Value = 7230
$ws.Range("D9").Value = 3713
$ws.Range("E9").Value = 1267

# Row 10
$ws.Range("B10").Value = 9051
$ws.Range("C10").Value = 6561
$ws.Range("D10").Value = 1917
$ws.Range("E10").Value = 573

# Row 13
$ws.Range("B13").Value = 5207
$ws.Range("C13").Value = 2704
$ws.Range("D13").Value = 1733
$ws.Range("E13").Value = 770

# Row 15
$ws.Range("B15").Value = 4936
$ws.Range("C15").Value = 2562
$ws.Range("D15").Value = 1903
$ws.Range("E15").Value = 471

# Row 16
$ws.Range("B16").Value = 3967
$ws.Range("C16").Value = 2345
$ws.Range("D16").Value = 1286
$ws.Range("E16").Value = 336

# Row 23
$ws.Range("B23").Value = 2852
$ws.Range("C23").Value = 2146
$ws.Range("D23").Value = 244
$ws.Range("E23").Value = 462

# Row 30
$ws.Range("B30").Value = 2308
$ws.Range("C30").Value = 936
$ws.Range("D30").Value = 1088

# Row 31
$ws.Range("B31").Value = 2231
$ws.Range("C31").Value = 1190
$ws.Range("D31").Value = 899
$ws.Range("E31").Value = 142

# Row 33
$ws.Range("B33").Value = 2207
$ws.Range("C33").Value = 1658
$ws.Range("D33").Value = 351
$ws.Range("E33").Value = 198
